$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.566.86"
$ws.Range("E2").Value = "  -5.39%  "
$ws.Range("D3").Value = "2.648.76"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.68"
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.56"
$ws.Range("E6").Value = "  -4.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.586"
$ws.Range("E7").Value = "  -2.75%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.566"
$ws.Range("E9").Value = "  -3.45%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.43"
$ws.Range("E10").Value = "  -5.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0829"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.91"
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("D13").Value = "3.055.18"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D15").Value = "2.665.86"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.907"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.81"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "44.615.35"
$ws.Range("E18").Value = "  -5.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").Value = "0.0₃0992"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.43"
$ws.Range("E21").Value = "  -4.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.07"
$ws.Range("E22").Value = "  +2.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "274.23"
$ws.Range("E23").Value = "  -0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.27"
$ws.Range("E24").Value = "  +4.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.98"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.34"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.35"
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.24"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.23"
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.08"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.28"
$ws.Range("E33").Value = "  +3.60%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "153.08"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.80"
$ws.Range("E35").Value = "  -2.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0822"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.118"
$ws.Range("E37").Value = "  -7.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.122"
$ws.Range("E38").Value = "  -0.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.27"
$ws.Range("E39").Value = "  +7.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.73"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.53"
$ws.Range("E41").Value = "  -3.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0318"
$ws.Range("E42").Value = "  -3.95%  "
$ws.Range("D43").Value = "2.124.68"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.88"
$ws.Range("E44").Value = "  -6.79%  "
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.05"
$ws.Range("E46").Value = "  -5.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.29"
$ws.Range("D48").Value = "2.912.98"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.77"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.194"
$ws.Range("E51").Value = "  -2.88%  "
